$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete rows 6 and 7 (target cluster "ECs" combos no longer present)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Update remaining rows 2-5 with refreshed TPM-derived NATMI statistics
$ws.Range("D2").Value = "FAPs"
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.294804333333333
$ws.Range("N2").Value = 6.884412999999999
$ws.Range("O2").Value = 0.6595351916216082
$ws.Range("P2").Value = 0.6595351916216082
$ws.Range("Q2").Value = 21.01342995829044
$ws.Range("R2").Value = 189.120869624614
$ws.Range("S2").Value = 0.6393347445319258
$ws.Range("T2").Value = 0.6393347445319258
$ws.Range("D3").Value = "MuSCs"
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("M3").Value = 1.184622333333333
$ws.Range("N3").Value = 3.553867
$ws.Range("O3").Value = 0.3404648083783919
$ws.Range("P3").Value = 0.3404648083783919
$ws.Range("Q3").Value = 10.84753853169178
$ws.Range("R3").Value = 97.62784678522601
$ws.Range("S3").Value = 0.3300369473106046
$ws.Range("T3").Value = 0.3300369473106046
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.2893236666666667
$ws.Range("H4").Value = 0.867971
$ws.Range("I4").Value = 0.03062830815746963
$ws.Range("J4").Value = 0.03062830815746962
$ws.Range("M4").Value = 2.294804333333333
$ws.Range("N4").Value = 6.884412999999999
$ws.Range("O4").Value = 0.6595351916216082
$ws.Range("P4").Value = 0.6595351916216082
$ws.Range("Q4").Value = 0.6639412040025555
$ws.Range("R4").Value = 5.975470836023
$ws.Range("S4").Value = 0.0202004470896824
$ws.Range("T4").Value = 0.02020044708968239
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.184622333333333
$ws.Range("N5").Value = 3.553867
$ws.Range("O5").Value = 0.3404648083783919
$ws.Range("P5").Value = 0.3404648083783919
$ws.Range("Q5").Value = 0.3427392770952222
$ws.Range("R5").Value = 3.084653493857
$ws.Range("S5").Value = 0.01042786106778723
$ws.Range("T5").Value = 0.01042786106778723

Write-Host "done"
